$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

$ws.Cells.Item($row, 1).Value = "'02/22/2026"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 9169.01
$ws.Cells.Item($row, 3).Value = 0.2538023500717298
$ws.Cells.Item($row, 4).Value = 0.7461976499282702
$ws.Cells.Item($row, 5).Value = -346.48
$ws.Cells.Item($row, 6).Value = -37.06
$ws.Cells.Item($row, 7).Value = -24149.86
$ws.Cells.Item($row, 8).Value = -77.92
$ws.Cells.Item($row, 9).Value = -1125.73
$ws.Cells.Item($row, 10).Value = -32.6
$ws.Cells.Item($row, 11).Value = -25275.59
$ws.Cells.Item($row, 12).Value = -73.38
